$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving the original
# "no explicit style" formatting (avoids Excel auto-converting
# numeric-looking strings like "27.76" or "597.24" into real numbers,
# and avoids leaving a stray NumberFormat/quote-prefix style behind).
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '75.975.66'
Set-TextValue $ws.Range('E2') '  +1.45%  '
Set-TextValue $ws.Range('D3') '2.912.93'
Set-TextValue $ws.Range('E3') '  +3.54%  '
Set-TextValue $ws.Range('E4') '  -0.03%  '
Set-TextValue $ws.Range('D5') '201.91'
Set-TextValue $ws.Range('E5') '  +7.47%  '
Set-TextValue $ws.Range('D6') '597.24'
Set-TextValue $ws.Range('E6') '  +0.49%  '
Set-TextValue $ws.Range('E7') '  -0.07%  '
Set-TextValue $ws.Range('E8') '  +0.10%  '
Set-TextValue $ws.Range('D9') '0.196'
Set-TextValue $ws.Range('E9') '  +2.67%  '
Set-TextValue $ws.Range('D10') '2.911.79'
Set-TextValue $ws.Range('E10') '  +3.53%  '
Set-TextValue $ws.Range('D11') '0.434'
Set-TextValue $ws.Range('E11') '  +16.84%  '
Set-TextValue $ws.Range('E12') '  +0.47%  '
Set-TextValue $ws.Range('E13') '  +0.43%  '
Set-TextValue $ws.Range('D14') '3.447.30'
Set-TextValue $ws.Range('E14') '  +3.47%  '
Set-TextValue $ws.Range('D15') '75.796.93'
Set-TextValue $ws.Range('E15') '  +1.21%  '
Set-TextValue $ws.Range('D16') '27.76'
Set-TextValue $ws.Range('E16') '  +3.51%  '
Set-TextValue $ws.Range('E17') '  +1.22%  '
Set-TextValue $ws.Range('D18') '2.913.94'
Set-TextValue $ws.Range('E18') '  +3.33%  '
Set-TextValue $ws.Range('D19') '12.87'
Set-TextValue $ws.Range('E19') '  +4.86%  '
Set-TextValue $ws.Range('D20') '8.73'
Set-TextValue $ws.Range('E20') '  -2.43%  '
Set-TextValue $ws.Range('D21') '371.83'
Set-TextValue $ws.Range('E21') '  -1.25%  '
Set-TextValue $ws.Range('D22') '2.32'
Set-TextValue $ws.Range('E22') '  +3.92%  '
Set-TextValue $ws.Range('D23') '4.29'
Set-TextValue $ws.Range('E23') '  +5.57%  '
Set-TextValue $ws.Range('B24') 'Litecoin'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D24') '71.08'
Set-TextValue $ws.Range('E24') '  +0.61%  '
Set-TextValue $ws.Range('B25') 'Dai'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D25') '0.999'
Set-TextValue $ws.Range('E25') '  -0.10%  '
Set-TextValue $ws.Range('D26') '3.053.75'
Set-TextValue $ws.Range('E26') '  +3.11%  '
Set-TextValue $ws.Range('D27') '4.23'
Set-TextValue $ws.Range('E27') '  +1.77%  '
Set-TextValue $ws.Range('D28') '9.68'
Set-TextValue $ws.Range('E28') '  +0.00%  '
Set-TextValue $ws.Range('E29') '  +2.82%  '
Set-TextValue $ws.Range('D30') '1.00'
Set-TextValue $ws.Range('E30') '  +0.25%  '
Set-TextValue $ws.Range('E31') '  +0.08%  '
Set-TextValue $ws.Range('D32') '501.33'
Set-TextValue $ws.Range('E32') '  -2.15%  '
Set-TextValue $ws.Range('D33') '7.69'
Set-TextValue $ws.Range('E33') '  -0.18%  '
Set-TextValue $ws.Range('E34') '  +1.89%  '
Set-TextValue $ws.Range('E35') '  -0.03%  '
Set-TextValue $ws.Range('D36') '165.68'
Set-TextValue $ws.Range('E36') '  +2.22%  '
Set-TextValue $ws.Range('E37') '  +1.22%  '
Set-TextValue $ws.Range('D38') '19.60'
Set-TextValue $ws.Range('E38') '  +1.16%  '
Set-TextValue $ws.Range('D39') '0.105'
Set-TextValue $ws.Range('E39') '  +23.95%  '
Set-TextValue $ws.Range('E40') '  -4.99%  '
Set-TextValue $ws.Range('E41') '  +0.00%  '
Set-TextValue $ws.Range('D42') '181.03'
Set-TextValue $ws.Range('E42') '  -3.28%  '
Set-TextValue $ws.Range('D43') '0.352'
Set-TextValue $ws.Range('E43') '  +4.25%  '
Set-TextValue $ws.Range('D44') '4.98'
Set-TextValue $ws.Range('E44') '  -0.41%  '
Set-TextValue $ws.Range('D45') '1.64'
Set-TextValue $ws.Range('E45') '  -1.13%  '
Set-TextValue $ws.Range('D46') '40.03'
Set-TextValue $ws.Range('E46') '  +0.05%  '
Set-TextValue $ws.Range('D48') '2.34'
Set-TextValue $ws.Range('E48') '  +0.67%  '
Set-TextValue $ws.Range('E50') '  +0.40%  '
Set-TextValue $ws.Range('D51') '0.655'
Set-TextValue $ws.Range('E51') '  +2.86%  '
